$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newA2 = @'
<Policy PolicyId="mli-interface" RuleCombiningAlgId="urn:oasis:names:tc:xacml:3.0:rule-combining-algorithm:permit-overrides" Version="1.0">
        <Target>
            <AnyOf>
                <AllOf>
                    <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                        <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">VIO</AttributeValue>
                        <AttributeDesignator AttributeId="http://authz-interop.org/AAA/xacml/subject/subject-role" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                    </Match>
                </AllOf>
            </AnyOf>
        </Target>
        <Rule Effect="Permit" RuleId="mli:replanning:vr-it-operations">
            <Target>
                <AnyOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">VR-IT</AttributeValue>
                            <AttributeDesignator AttributeId="http://authz-interop.org/AAA/xacml/resource/resource-type" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                </AnyOf>
                <AnyOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:ReplanningVI:Add-VR-IT</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:ReplanningVI:Modify-VR-IT</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:ReplanningVI:Delete-VR-IT</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                </AnyOf>
            </Target>
        </Rule>
    </Policy>
'@

$newC2 = @'
{
    "datalog_subjects": "Subject(S).",
    "datalog_objects": "Resource(R).",
    "datalog_relationships": "has_role(S, \"VIO\") :- Subject(S).\nhas_type(R, \"VR-IT\") :- Resource(R).",
    "datalog_actions": "can_perform(S, R, \"MLI:ReplanningVI:Add-VR-IT\") :- Subject(S), Resource(R), has_role(S, \"VIO\"), has_type(R, \"VR-IT\").\ncan_perform(S, R, \"MLI:ReplanningVI:Modify-VR-IT\") :- Subject(S), Resource(R), has_role(S, \"VIO\"), has_type(R, \"VR-IT\").\ncan_perform(S, R, \"MLI:ReplanningVI:Delete-VR-IT\") :- Subject(S), Resource(R), has_role(S, \"VIO\"), has_type(R, \"VR-IT\")."
}
'@

$newC3 = @'
{
    "datalog_subjects": "VIO(S).",
    "datalog_objects": "VR(R), VR_RP_Info(R), VR_State_Info(R), VR_Power_Info(R), VR_Status_Info(R), VR_Info(R).",
    "datalog_relationships": "",
    "datalog_actions": "can_perform(S, R, Decommission_VR_IT) :- VIO(S), VR(R). can_perform(S, R, Operate_VR) :- VIO(S), VR(R). can_perform(S, R, Add_VirtualNetworkIf) :- VIO(S), VR(R). can_perform(S, R, Remove_VirtualNetworkIf) :- VIO(S), VR(R). can_perform(S, R, Create_StorageImage) :- VIO(S), VR(R). can_perform(S, R, Remove_StorageImage) :- VIO(S), VR(R). can_perform(S, R, Get_Available_VR_Pool_IT) :- VIO(S), VR_RP_Info(R). can_perform(S, R, Monitor_VR_Info) :- VIO(S), VR_State_Info(R), VR_Power_Info(R), VR_Status_Info(R). can_perform(S, R, Subscribe_VR_Monitoring) :- VIO(S), VR_Info(R). can_perform(S, R, Unsubscribe_VR_Monitoring) :- VIO(S), VR_Info(R)."
}
'@

$ws.Range("A2").Value = $newA2
$ws.Range("C2").Value = $newC2
$ws.Range("C3").Value = $newC3
